# df_gwent_leaders.xlsx fixes + patch 7.1 data additions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# -----------------------------------------------------------------
# 1) Tooltip (column K) text corrections / reworked numbers
# -----------------------------------------------------------------
$ws.Range("K20").Value = "Order: Play a Skellige unit with provision 9 or less from your graveyard.  This ability adds 12 provisions to your deck's provisions limit. "

$ws.Range("K41").Value = "Order: Destroy an allied unit, then Spawn an Ekimmara in its row and boost it by the destroyed unit's power. Charge: 3.  This ability adds 14 provisions to your deck's provision limit. "
$ws.Range("K42").Value = "Order: Lock an enemy unit, then Spawn its base copy in the opposite row and boost it by 2.  This ability adds 14 provisions to your deck's provision limit. "
$ws.Range("K43").Value = "Order: Boost an allied unit by 1. Charge: 3 Once all Charges are used up, Spawn and play Lyrian Scytheman.  This ability adds 14 provisions to your deck's provision limit. "
$ws.Range("K44").Value = "Order: Spawn and play Tempering. All Dwarves in your starting deck get 1 Armor.  This ability adds 16 provisions to your deck's provision limit. "
$ws.Range("K45").Value = "Order: Move a Skellige unit from your deck to your graveyard, then damage an enemy unit by its power.  This ability adds 16 provisions to your deck's provision limit. "
$ws.Range("K46").Value = "Order: Gain 3 Coins. At the beginning of the round, refresh this ability. Your Hoards require 2 less Coins to trigger.  This ability adds 14 provisions to your deck's provision limit. "

# -----------------------------------------------------------------
# 2) Provision (column E) value corrections
# -----------------------------------------------------------------
$ws.Range("E20").Value = 12
$ws.Range("E41").Value = 14
$ws.Range("E42").Value = 14
$ws.Range("E44").Value = 16
$ws.Range("E45").Value = 16
$ws.Range("E46").Value = 14

# -----------------------------------------------------------------
# 3) Row-id (column A) renumbering, rows 4-46
# -----------------------------------------------------------------
$idUpdates = @{
    4  = 114
    5  = 340
    6  = 375
    7  = 376
    8  = 377
    9  = 379
    10 = 380
    11 = 381
    12 = 382
    13 = 383
    14 = 384
    15 = 385
    16 = 386
    17 = 441
    18 = 448
    19 = 449
    20 = 450
    21 = 483
    22 = 530
    23 = 531
    24 = 551
    25 = 569
    26 = 570
    27 = 571
    28 = 572
    29 = 573
    30 = 581
    31 = 597
    32 = 598
    33 = 599
    34 = 600
    35 = 601
    36 = 711
    37 = 735
    38 = 754
    39 = 780
    40 = 797
    41 = 957
    42 = 958
    43 = 959
    44 = 960
    45 = 961
    46 = 962
}

foreach ($row in $idUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value = $idUpdates[$row]
}
